$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 50005532
$ws.Range("I18").Value = 50005532
$ws.Range("K18").Value = 50005532
$ws.Range("M18").Value = -50005248
$ws.Range("H40").Value = 3698.9
$ws.Range("I40").Value = 3328.5715
$ws.Range("J40").Value = 4563
$ws.Range("K40").Value = 3328.5715
$ws.Range("L40").Value = 4563
$ws.Range("M40").Value = -3153.5715
$ws.Range("N40").Value = -4913
$ws.Range("H103").Value = 1580.4546
$ws.Range("I103").Value = 789.5
$ws.Range("K103").Value = 2368.5
$ws.Range("M103").Value = -1782.5
$ws.Range("H112").Value = 4846.5293
$ws.Range("I112").Value = 971
$ws.Range("K112").Value = 2913
$ws.Range("M112").Value = -1805
$ws.Range("H113").Value = 100006630
$ws.Range("J113").Value = 136371630
$ws.Range("L113").Value = 136371630
$ws.Range("N113").Value = -136378138
$ws.Range("H125").Value = 83335336
$ws.Range("I125").Value = 500000000
$ws.Range("J125").Value = 2400.2
$ws.Range("K125").Value = 4500000000
$ws.Range("L125").Value = 21601.8
$ws.Range("M125").Value = -4499997540
$ws.Range("N125").Value = -26521.8
$ws.Range("H132").Value = 1296.3864
$ws.Range("I132").Value = 1251.8536
$ws.Range("K132").Value = 3755.5608
$ws.Range("M132").Value = -1225.5608
$ws.Range("H137").Value = 3098.2
$ws.Range("I137").Value = 3037.8462
$ws.Range("J137").Value = 3490.5
$ws.Range("K137").Value = 9113.5386
$ws.Range("L137").Value = 10471.5
$ws.Range("M137").Value = -6563.5386
$ws.Range("N137").Value = -15571.5
$ws.Range("H138").Value = 4066.6824
$ws.Range("I138").Value = 1185.52
$ws.Range("K138").Value = 3556.56
$ws.Range("M138").Value = 1583.44

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 4210
$ws.Range("I10").Value = 4210
$ws.Range("K10").Value = 4210
$ws.Range("M10").Value = -4040
$ws.Range("H32").Value = 1609119.8
$ws.Range("I32").Value = 1694822
$ws.Range("K32").Value = 1694822
$ws.Range("M32").Value = -1694535
$ws.Range("H45").Value = 6820.7
$ws.Range("I45").Value = 4601
$ws.Range("K45").Value = 4601
$ws.Range("M45").Value = -4224
$ws.Range("H61").Value = 5220.7964
$ws.Range("I61").Value = 3033.348
$ws.Range("J61").Value = 12961
$ws.Range("K61").Value = 3033.348
$ws.Range("L61").Value = 12961
$ws.Range("M61").Value = -2821.348
$ws.Range("N61").Value = -13385
$ws.Range("H129").Value = 82049.836
$ws.Range("J129").Value = 82049.836
$ws.Range("L129").Value = 82049.836
$ws.Range("N129").Value = -92049.836
$ws.Range("H132").Value = 6185.879
$ws.Range("I132").Value = 1853.9166
$ws.Range("J132").Value = 8661.286
$ws.Range("K132").Value = 5561.7498
$ws.Range("L132").Value = 25983.858
$ws.Range("M132").Value = -3031.7498
$ws.Range("N132").Value = -31043.858
$ws.Range("H136").Value = 5220.7964
$ws.Range("I136").Value = 3033.348
$ws.Range("J136").Value = 12961
$ws.Range("K136").Value = 9100.044
$ws.Range("L136").Value = 38883
$ws.Range("M136").Value = -6550.044
$ws.Range("N136").Value = -43983

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 104172830
$ws.Range("I86").Value = 50004800
$ws.Range("J86").Value = 142864290
$ws.Range("K86").Value = 50004800
$ws.Range("L86").Value = 142864290
$ws.Range("M86").Value = -50003677
$ws.Range("N86").Value = -142866536
$ws.Range("H89").Value = 104172830
$ws.Range("I89").Value = 50004800
$ws.Range("J89").Value = 142864290
$ws.Range("K89").Value = 250024000
$ws.Range("L89").Value = 714321450
$ws.Range("M89").Value = -250018384
$ws.Range("N89").Value = -714332682
$ws.Range("H94").Value = 2522.7896
$ws.Range("I94").Value = 742.8333
$ws.Range("J94").Value = 5574.143
$ws.Range("K94").Value = 742.8333
$ws.Range("L94").Value = 5574.143
$ws.Range("M94").Value = -291.8333
$ws.Range("N94").Value = -6476.143
$ws.Range("H99").Value = 45457044
$ws.Range("I99").Value = 5000
$ws.Range("K99").Value = 5000
$ws.Range("M99").Value = -3502
$ws.Range("H107").Value = 75003610
$ws.Range("J107").Value = 5559.222
$ws.Range("L107").Value = 5559.222
$ws.Range("N107").Value = -9399.222
$ws.Range("H134").Value = 7752.727
$ws.Range("I134").Value = 3211.4443
$ws.Range("J134").Value = 9455.708000000001
$ws.Range("K134").Value = 9634.332900000001
$ws.Range("L134").Value = 28367.124
$ws.Range("M134").Value = -7099.332900000001
$ws.Range("N134").Value = -33437.124

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 515.75
$ws.Range("I5").Value = 519.1429000000001
$ws.Range("J5").Value = 492
$ws.Range("K5").Value = 519.1429000000001
$ws.Range("L5").Value = 492
$ws.Range("M5").Value = -407.1429000000001
$ws.Range("N5").Value = -716
$ws.Range("H22").Value = 107.666664
$ws.Range("I22").Value = 67.28570999999999
$ws.Range("K22").Value = 67.28570999999999
$ws.Range("M22").Value = 282.71429
$ws.Range("H31").Value = 5273.8203
$ws.Range("I31").Value = 2184.5881
$ws.Range("J31").Value = 11109.037
$ws.Range("K31").Value = 2184.5881
$ws.Range("L31").Value = 11109.037
$ws.Range("M31").Value = -1889.5881
$ws.Range("N31").Value = -11699.037
$ws.Range("H34").Value = 5273.8203
$ws.Range("I34").Value = 2184.5881
$ws.Range("J34").Value = 11109.037
$ws.Range("K34").Value = 2184.5881
$ws.Range("L34").Value = 11109.037
$ws.Range("M34").Value = -1982.5881
$ws.Range("N34").Value = -11513.037
$ws.Range("H58").Value = 7466790.5
$ws.Range("I58").Value = 10639658
$ws.Range("J58").Value = 10551.85
$ws.Range("K58").Value = 10639658
$ws.Range("L58").Value = 10551.85
$ws.Range("M58").Value = -10639455
$ws.Range("N58").Value = -10957.85
$ws.Range("H62").Value = 9986.625
$ws.Range("I62").Value = 9984
$ws.Range("K62").Value = 9984
$ws.Range("M62").Value = -9360
$ws.Range("H65").Value = 9986.625
$ws.Range("I65").Value = 9984
$ws.Range("K65").Value = 49920
$ws.Range("M65").Value = -46800
$ws.Range("H127").Value = 50318.8
$ws.Range("J127").Value = 50318.8
$ws.Range("L127").Value = 50318.8
$ws.Range("N127").Value = -60238.8
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H134").Value = 4190.29
$ws.Range("I134").Value = 1740.8085
$ws.Range("K134").Value = 5222.4255
$ws.Range("M134").Value = -2687.4255
$ws.Range("H136").Value = 7466790.5
$ws.Range("I136").Value = 10639658
$ws.Range("J136").Value = 10551.85
$ws.Range("K136").Value = 31918974
$ws.Range("L136").Value = 31655.55
$ws.Range("M136").Value = -31916424
$ws.Range("N136").Value = -36755.55

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1292502.9
$ws.Range("I5").Value = 2667286.8
$ws.Range("J5").Value = 3643.125
$ws.Range("K5").Value = 8001860.399999999
$ws.Range("L5").Value = 10929.375
$ws.Range("M5").Value = -8001748.399999999
$ws.Range("N5").Value = -11153.375
$ws.Range("H92").Value = 6412048
$ws.Range("I92").Value = 1805
$ws.Range("J92").Value = 7694096.5
$ws.Range("K92").Value = 5415
$ws.Range("L92").Value = 23082289.5
$ws.Range("M92").Value = -4167
$ws.Range("N92").Value = -23084785.5
$ws.Range("H98").Value = 2488
$ws.Range("J98").Value = 2830
$ws.Range("L98").Value = 8490
$ws.Range("N98").Value = -11486
$ws.Range("H103").Value = 1267.3846
$ws.Range("I103").Value = 573.75
$ws.Range("K103").Value = 1721.25
$ws.Range("M103").Value = -842.25
$ws.Range("H129").Value = 72928.64
$ws.Range("J129").Value = 84828.414
$ws.Range("L129").Value = 254485.242
$ws.Range("N129").Value = -264485.242
$ws.Range("H135").Value = 1292502.9
$ws.Range("I135").Value = 2667286.8
$ws.Range("J135").Value = 3643.125
$ws.Range("K135").Value = 24005581.2
$ws.Range("L135").Value = 32788.125
$ws.Range("M135").Value = -24003046.2
$ws.Range("N135").Value = -37858.125

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 10000
$ws.Range("I44").Value = 10000
$ws.Range("K44").Value = 10000
$ws.Range("M44").Value = -9404
$ws.Range("H141").Value = 27831.584
$ws.Range("J141").Value = 32110.666
$ws.Range("L141").Value = 32110.666
$ws.Range("N141").Value = -42470.666

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1223.4445
$ws.Range("I22").Value = 546.4545000000001
$ws.Range("K22").Value = 546.4545000000001
$ws.Range("M22").Value = -251.4545000000001
$ws.Range("H27").Value = 1223.4445
$ws.Range("I27").Value = 546.4545000000001
$ws.Range("K27").Value = 546.4545000000001
$ws.Range("M27").Value = -439.4545000000001
$ws.Range("H68").Value = 7667.3335
$ws.Range("I68").Value = 3002
$ws.Range("K68").Value = 3002
$ws.Range("M68").Value = -2253
$ws.Range("H71").Value = 7667.3335
$ws.Range("I71").Value = 3002
$ws.Range("K71").Value = 15010
$ws.Range("M71").Value = -11266
$ws.Range("H136").Value = 6481.7383
$ws.Range("I136").Value = 968.25
$ws.Range("K136").Value = 2904.75
$ws.Range("M136").Value = -354.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H116").Value = 55784.5
$ws.Range("J116").Value = 55784.5
$ws.Range("L116").Value = 55784.5
$ws.Range("N116").Value = -64962.5
$ws.Range("H132").Value = 12205685
$ws.Range("I132").Value = 16133680
$ws.Range("J132").Value = 28899.9
$ws.Range("K132").Value = 48401040
$ws.Range("L132").Value = 86699.70000000001
$ws.Range("M132").Value = -48398510
$ws.Range("N132").Value = -91759.70000000001
